$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove all hyperlinks (Username/Password/Email/PhoneNo columns had
#     mailto: hyperlinks on the Email column). The collection mutates as
#     items are removed, so sweep repeatedly until it is empty. ---
$guard = 0
while ($ws.Hyperlinks.Count -gt 0 -and $guard -lt 10) {
    $guard = $guard + 1
    $snapshot = @($ws.Hyperlinks)
    foreach ($link in $snapshot) {
        $link.Delete()
    }
}

# --- Drop the now-unneeded Username/Password/Email/PhoneNo columns (F:I) ---
$ws.Range("F1:I1").EntireColumn.Delete()

# --- Resize the Staff ID column (mirrors the old "Password" column width) ---
$ws.Columns("A").ColumnWidth = 15

# --- Renumber the staff IDs from department-prefixed codes to S0xx codes ---
$ws.Range("A2").Value = "S001"
$ws.Range("A3").Value = "S002"
$ws.Range("A4").Value = "S003"
$ws.Range("A5").Value = "S004"

# --- The Hyperlink cell style is no longer used anywhere in the sheet ---
$wb.Styles("Hyperlink").Delete()

# --- Restore the active selection ---
[void]$ws.Range("C12").Select()
